# Add two new rows (10 and 11) to the translation table on the active sheet,
# mirroring the existing "filename row" + "data row" layout used throughout
# the sheet (see rows 2/3, 4, 5, ... 8, 9).
#
# Row 10 reuses the formatting of row 3 (the "continuation" style: s=6/7)
# Row 11 reuses the formatting of row 9 (the "main" style: s=4/5)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy cell formatting from existing rows into the new rows -------------
$ws.Range("A3:E3").Copy()
$ws.Range("A10:E10").PasteSpecial(-4122)

$ws.Range("A9:E9").Copy()
$ws.Range("A11:E11").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Row heights -------------------------------------------------------
$ws.Rows.Item(10).RowHeight = 43.2
$ws.Rows.Item(11).RowHeight = 43.2

# --- Set new cell values -------------------------------------------------
# The order of assignment below matches the order new strings were added
# to the shared-string table in the committed workbook (filename, English
# line, second filename, Russian translation, then the "converted" line).
$ws.Range("A10").Value2 = "SCRIPT/G01P03A/um2505.ssb"
$ws.Range("C11").Value2 = " The grand master of all things\nbad?[K] ...What?"
$ws.Range("A11").Value2 = "SCRIPT/G01P03A/us0106.ssb"
$ws.Range("D11").Value2 = " Гранд мастер всего самого\nплохого?[K] ...Что?"
$ws.Range("E11").Value2 = " Ãñàîä íàòóåñ âòåãï òàíïãï\nðìïöïãï?[K] ...Œóï?"
$ws.Range("B11").Value2 = 76

# --- Update selection to match the committed workbook state --------------
$ws.Range("D8").Select()
